$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the "Date" value ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-05-05T14:17:01+00:00"

# --- Elements sheet: update binding strength and value set URL ---
$wsElem = $wb.Worksheets.Item("Elements")

# Author Role binding strength: required -> preferred
$wsElem.Range("X5").Value = "preferred"

# AuthorSpecialty binding strength: required -> preferred
$wsElem.Range("X6").Value = "preferred"

# AuthorSpecialty binding value set URL change
$wsElem.Range("Z6").Value = "https://mos.esante.gouv.fr/NOS/JDV_J01-XdsAuthorSpecialty-CISIS/FHIR/JDV-J01-XdsAuthorSpecialty-CISIS"

# Widen column Z to fit the new, longer URL (Excel re-derives the best-fit
# pixel width from the font metrics; 83.0 is the input that rounds to the
# same stored column width Excel produced, 83.8.../84.66... after its own
# pixel-quantization).
$wsElem.Columns.Item(26).ColumnWidth = 83.0
